$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "A Software Engineer designs, develops, tests, and maintains software applications or systems. They work with various programming languages, frameworks, and tools to ensure the software functions correctly, is scalable, and is user-friendly. They are also responsible for debugging and improving existing code, collaborating with other developers, and following industry best practices and coding standards."

$ws.Range("B3").Value = "A Project Manager oversees the planning, execution, and completion of projects within a set timeframe and budget. They define project scope, manage resources, set milestones, and ensure effective communication among teams. They also identify potential risks and develop mitigation strategies, ensuring the project stays within the agreed-upon parameters. Project Managers often act as the primary point of contact between stakeholders and the project team."

$ws.Range("B4").Value = "An HR Specialist is responsible for managing the human resources functions of an organization, including recruitment, employee relations, performance management, training and development, and compliance with labor laws. They ensure that the organization has the right talent to meet its goals while fostering a positive workplace culture. HR specialists also handle administrative tasks such as payroll processing, benefits administration, and maintaining employee records."

$ws.Range("B5").Value = "An Accountant is responsible for managing and maintaining financial records, ensuring compliance with regulatory requirements, and preparing financial statements. They handle budgeting, forecasting, tax preparation, and financial reporting. Accountants work to ensure that a company’s financial activities are conducted in accordance with accounting standards and help provide accurate data to support decision-making. They may also assist in auditing processes and provide recommendations for financial improvements."

$ws.Range("B6").Value = "A Marketing Manager develops and implements marketing strategies to promote products or services and drive business growth. They oversee market research, identify customer needs, and manage promotional campaigns across various channels, including digital marketing, advertising, public relations, and events. The Marketing Manager works closely with sales, product development, and creative teams to ensure alignment with business objectives and effectively reach target audiences."

$ws.Range("B7").Value = "Customer Support professionals assist customers by providing answers to inquiries, resolving issues, and guiding them through product or service usage. They are responsible for ensuring a positive customer experience by addressing concerns, troubleshooting problems, and offering solutions. Customer Support agents often handle customer feedback, provide product knowledge, and escalate more complex issues to other departments as needed. They play a critical role in customer retention and satisfaction."

$ws.Range("B8").Value = "An Operations Manager ensures that the day-to-day operations of an organization run efficiently and effectively. They oversee operational processes, manage logistics, coordinate with other departments, and improve operational systems to optimize productivity and reduce costs. They also manage the supply chain, handle resource allocation, and ensure that quality standards are met. The Operations Manager works to streamline workflows, improve processes, and ensure the organization’s goals are achieved on time."

$ws.Range("B9").Value = "A Product Designer focuses on creating user-centered designs for products, ensuring that both functionality and aesthetics are balanced. They conduct user research, create wireframes, prototypes, and visual designs, and work closely with engineering teams to ensure that the design is feasible and meets user needs. Product Designers also test prototypes, gather feedback, and iterate on designs to enhance user experiences. They strive to make products intuitive, visually appealing, and accessible."

$ws.Range("B10").Value = "A Business Analyst analyzes business processes, identifies opportunities for improvement, and provides data-driven insights to support decision-making. They work with stakeholders to gather requirements, analyze data, and develop solutions that help streamline operations, reduce costs, or enhance customer satisfaction. Business Analysts also help create strategies for business growth, define key performance indicators (KPIs), and monitor progress toward business objectives."

$ws.Range("B11").Value = "A Data Scientist analyzes large sets of complex data to uncover patterns, trends, and insights that can inform business decisions. They apply statistical models, machine learning algorithms, and advanced analytics techniques to interpret data and make predictions. Data Scientists work with cross-functional teams to define business challenges, clean and preprocess data, build models, and communicate insights. They also create automated systems for data analysis and ensure that the models are scalable and efficient."

$ws.Range("B12").Value = "A Data Analyst interprets data and provides actionable insights that help businesses make informed decisions. They collect, process, and analyze data from various sources, using statistical methods and data visualization techniques to identify trends, patterns, and correlations. Data Analysts often work with stakeholders to define data needs, generate reports, and present findings that support decision-making processes. They play a key role in improving business operations and strategies by providing data-driven recommendations."

$ws.Range("B13").Value = "A Machine Learning Engineer builds and deploys machine learning models that enable systems to learn from data and improve over time. They design algorithms, train models, and develop systems that can automate processes or make predictions. Machine Learning Engineers work with large datasets, clean and prepare data, and optimize models for efficiency. They collaborate with data scientists and software engineers to integrate machine learning models into applications or systems, and ensure that the models perform well in production environments."

$ws.Range("B14").Value = "A Sales Executive is responsible for driving sales and revenue growth by identifying business opportunities, building relationships with clients, and closing deals. They work to understand customer needs, present products or services, negotiate terms, and close sales. Sales Executives often work with marketing teams to generate leads, develop sales strategies, and meet sales targets. They maintain relationships with existing customers, follow up on leads, and continuously identify new sales opportunities."

$ws.Range("B15").Value = "A Legal Advisor provides legal counsel to organizations or individuals, ensuring compliance with laws and regulations. They review contracts, negotiate terms, and provide advice on legal risks and business decisions. Legal Advisors help organizations navigate complex legal frameworks, handle disputes, and represent the company in legal matters. They ensure that the business operates within the boundaries of the law and may also assist with intellectual property rights, corporate governance, and regulatory compliance."

$wb.Save()
